$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new row 23 that mirrors row 22 (Organization) but for Person/PERSON
$ws.Range("B23").Value = "Person - Check participants list for NoAccess & Owner"
$ws.Range("C23").Value = "PERSON"
$ws.Range("D23").Value = $ws.Range("D22").Value2
$ws.Range("E23").Value = $ws.Range("E22").Value2
$ws.Range("F23").Value = $ws.Range("F22").Value2

# Copy formatting from row 22 to the new row 23
$ws.Range("B22:F22").Copy()
$ws.Range("B23:F23").PasteSpecial(-4122)
$ws.Rows.Item(23).RowHeight = 45

$ws.Range("F23").Select()
